$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("S2").Value = 0.5
$ws.Range("U2").Value = 0.5

$ws.Range("C3").Value = 0.5
$ws.Range("F3").Value = 1
$ws.Range("S3").Value = 0.5
$ws.Range("U3").Value = 0.5

$ws.Range("C4").Value = 0.5
$ws.Range("F4").Value = 0.5
$ws.Range("S4").Value = 0.5
$ws.Range("U4").Value = 0.5

$ws.Range("C5").Value = 0.5
$ws.Range("F5").Value = 0.7142857142857143
$ws.Range("S5").Value = 0.5
$ws.Range("U5").Value = 0.5

$ws.Range("C6").Value = 0.2754115523761866
$ws.Range("F6").Value = 1
$ws.Range("S6").Value = 0.17376534287144
$ws.Range("U6").Value = 0.2754115523761866
